$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "for Qualitative Research Group" -> "for the Qualitative Research Group"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Wrote Lisp code (mainly GUI) for Qualitative Research Group.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Wrote Lisp code (mainly GUI) for the Qualitative Research Group.", 2)

# ---------------------------------------------------------------------
# 2) "Between degrees, did a work abroad" -> "Between degrees, I did a work abroad"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Between degrees, did a work abroad", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Between degrees, I did a work abroad", 2)

# ---------------------------------------------------------------------
# 3) "highlight undergrad degrees too"
#    - "B.S. Physics and B.S. Biophysics" grows to the (unset) document
#      default size instead of the small 10pt used by the rest of the
#      résumé line (closest reachable approximation: bump the Font.Size
#      explicitly since the COM layer can't clear an explicit run size).
#    - " ... half-time research programmer" gets respelled as
#      "... ½ time research-programmer" with "½ time" and "-" losing
#      their explicit black-color run property (closest reachable
#      approximation: Font.Color -> Automatic).
# ---------------------------------------------------------------------
$deg = $d.Content
$deg.Find.Execute("B.S. Physics and B.S. Biophysics", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deg.Font.Size = 12

$tail = $d.Content
$tail.Find.Execute("half-time research programmer", $true, $false, $false, $false, $false, $true, 1, $false,
                    "½ time research-programmer", 2)

$halfTime = $d.Content
$halfTime.Find.Execute("½ time", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$halfTime.Font.Color = -16777216

$hyphenHost = $d.Content
$hyphenHost.Find.Execute("research-programmer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hyphenPos = $hyphenHost.Start + 8
$hyphen = $d.Range($hyphenPos, $hyphenPos + 1)
$hyphen.Font.Color = -16777216
